$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "This beautiful baby‑centered opera—based on a tender tale of a mother bird and chick—creates an immersive sensory world full of bright colors, gentle melodies, and interactive play designed for infants." "A delightful sensory experience, BambinO invites babies and caregivers into a cozy immersive opera: a charming mother-bird narrative filled with engaging colors, sounds, and opportunities for little ones to chirp and play through their imagination—all in a safe, age-appropriate setting."

Replace-Text "📍 Location: Topanga Library" "📍 Location: Child Development Institute, Canoga Park"

Replace-Text "122 N. Topanga Canyon Blvd, Topanga, CA 90290" "7260 Owensmouth Ave, Canoga Park, CA 91303"

Replace-Text "📅 Date: 2025‑06‑13" "📅 Date: 2025‑06‑20"

Replace-Text "🕘 Time: 11:00 AM" "🕘 Time: 10 AM & 2 PM"

Replace-Text "👶 Age Requirement: Babies 6–18 months" "👶 Age Requirement: Ages 6–18 mo"

Replace-Text "#BambinO #BabyOpera #InteractiveMusic #FreeEvent #LAOperaConnects #InfantImagination #ParentChild #MusicalPlay #LibraryEvent #ShitToDoWithKids #shittodowithkids #stdwkids #familyactivities #kidslosangeles" "#babies #infantopera #interactive #earlymusic #laopera #freeevent #westvalley #canogapark #bambino #ShitToDoWithKids #shittodowithkids #stdwkids #familyactivities #kidslosangeles"
